$wb = $excel.ActiveWorkbook

# "Generate Report for Handback": record a new handoff/handback round-trip
# for the "ebee5ea8-fce3-4e2e-8c60-50f1dae031a9" file (row 7) on both the
# zh-cn and de-de localization sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value = "2016-03-04 08:30:06"
$wsZhCn.Range("G7").Value = "2016-03-04 08:30:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value = "2016-03-04 08:30:17"
$wsDeDe.Range("G7").Value = "2016-03-04 08:31:20"
